# Applies the "gh-pages output generated at 456a3b4" update:
#   - Sheet "展览" (exhibitions): refresh "want-to-go" counters on several
#     existing rows, update one cover-image URL, and insert a brand-new
#     exhibition row (北京·地狱双ip同人ONLY展, 2024-09-22) ahead of the
#     2024-10-01 entries, shifting everything below it down by one row.
#   - Sheet "演出" (performances): refresh counters; one listing's lowest
#     price flips from a number to the literal text "不可售" (sold out /
#     not orderable).
#   - Sheet "本地生活" (local life): refresh one counter.
#   - Sheet "全部类型" (all types, an aggregate roll-up): refresh the same
#     counters mirrored from the three sheets above (this sheet is not
#     re-sorted/re-inserted with the new exhibition row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Simple counter refreshes on rows that are NOT affected by the insert
# (they sit above the insertion point, row 30).
$ws1.Cells.Item(4, 6).Value  = 149   # was 147
$ws1.Cells.Item(6, 6).Value  = 3736  # was 3728
$ws1.Cells.Item(8, 6).Value  = 2514  # was 2506
$ws1.Cells.Item(9, 6).Value  = 63    # was 60
$ws1.Cells.Item(10, 6).Value = 3008  # was 2993
$ws1.Cells.Item(16, 6).Value = 77    # was 73
$ws1.Cells.Item(17, 6).Value = 429   # was 428
$ws1.Cells.Item(19, 6).Value = 186   # was 185
$ws1.Cells.Item(20, 6).Value = 335   # was 333
$ws1.Cells.Item(21, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/diA04BbQ1724055935646.jpeg"
$ws1.Cells.Item(22, 6).Value = 318   # was 316
$ws1.Cells.Item(23, 6).Value = 633   # was 632
$ws1.Cells.Item(24, 6).Value = 1377  # was 1374

# Insert a brand-new row at position 30 (everything from the old row 30
# downward shifts to row+1). Copy formatting only from the row above so
# the new row's index cell (column A) keeps the bold/centered/bordered
# style used throughout column A.
$ws1.Rows.Item(29).Copy() | Out-Null
$ws1.Rows.Item(30).Insert() | Out-Null

$ws1.Cells.Item(30, 1).Value = 29
$ws1.Cells.Item(30, 2).NumberFormat = "@"
$ws1.Cells.Item(30, 2).Value = "2024-09-22"
$ws1.Cells.Item(30, 3).Value = "北京·地狱双ip同人ONLY展"
$ws1.Cells.Item(30, 4).Value = "双桥中路50号院 E50艺术园区"
$ws1.Cells.Item(30, 5).Value = "2024.09.22 10:30-09.22 16:00"
$ws1.Cells.Item(30, 6).Value = 0
$ws1.Cells.Item(30, 7).Value = 95
$ws1.Cells.Item(30, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90931"
$ws1.Cells.Item(30, 9).Value = "//i1.hdslb.com/bfs/openplatform/202408/c6ObwO4C1724055713128.jpeg"

# Counter refreshes on the rows that shifted down by one (now at their
# +1 positions). Rows whose counters did not change are left untouched.
$ws1.Cells.Item(31, 6).Value = 4135  # 北京·IDO动漫游戏嘉年华47th, was 4119
$ws1.Cells.Item(32, 6).Value = 3735  # 北京·第19届IJOY漫展xCGF游戏节, was 3710
$ws1.Cells.Item(33, 6).Value = 62    # 北京·第五人格同人only同人3.0, was 60
$ws1.Cells.Item(35, 6).Value = 1095  # 北京·第五人格only同人展, was 1094
$ws1.Cells.Item(36, 6).Value = 446   # 北京·ICOS国际动漫节×CGF中国游戏节04, was 445
$ws1.Cells.Item(38, 6).Value = 1302  # 北京·明日方舟only同人展, was 1301
$ws1.Cells.Item(39, 6).Value = 142   # 北京·第一届世界计划pjsk only同人展, was 141
$ws1.Cells.Item(41, 6).Value = 86    # 北京·Aw动漫游戏嘉年华9th, was 85
$ws1.Cells.Item(43, 6).Value = 55    # 北京·万游引力国潮动漫嘉年华S9, was 54
$ws1.Cells.Item(44, 6).Value = 48    # 北京·代号鸢only同人展, was 45

# ---------------------------------------------------------------------
# Sheet 2: 演出 (performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Cells.Item(8, 6).Value = 19   # was 17
$ws2.Cells.Item(9, 6).Value = 14   # was 13
$ws2.Cells.Item(10, 7).Value = "不可售"  # was numeric 180 (sold out)
$ws2.Cells.Item(15, 6).Value = 42  # was 41
$ws2.Cells.Item(16, 6).Value = 197 # was 196

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Cells.Item(4, 6).Value = 2217  # was 2210

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (all types - aggregate roll-up, no row insert here)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(7, 6).Value  = 149   # was 147
$ws4.Cells.Item(11, 6).Value = 3736  # was 3728
$ws4.Cells.Item(13, 6).Value = 2514  # was 2506
$ws4.Cells.Item(14, 6).Value = 63    # was 60
$ws4.Cells.Item(15, 6).Value = 3008  # was 2993
$ws4.Cells.Item(20, 6).Value = 77    # was 73
$ws4.Cells.Item(21, 6).Value = 429   # was 428
$ws4.Cells.Item(23, 6).Value = 335   # was 333
$ws4.Cells.Item(24, 6).Value = 318   # was 316
$ws4.Cells.Item(25, 6).Value = 633   # was 632
$ws4.Cells.Item(26, 6).Value = 1377  # was 1374
$ws4.Cells.Item(32, 6).Value = 19    # was 17
$ws4.Cells.Item(33, 6).Value = 4135  # was 4119
$ws4.Cells.Item(34, 6).Value = 3735  # was 3710
$ws4.Cells.Item(35, 6).Value = 62    # was 60
$ws4.Cells.Item(38, 6).Value = 446   # was 445
$ws4.Cells.Item(42, 6).Value = 42    # was 41
$ws4.Cells.Item(43, 6).Value = 1302  # was 1301
$ws4.Cells.Item(44, 6).Value = 142   # was 141
$ws4.Cells.Item(45, 6).Value = 86    # was 85
$ws4.Cells.Item(47, 6).Value = 55    # was 54
$ws4.Cells.Item(48, 6).Value = 48    # was 45
$ws4.Cells.Item(49, 6).Value = 197   # was 196

Write-Output "edit complete"
